# Append the weekly bitcoin-buy record for 2025-12-21 as row 62,
# matching the existing layout: col A is the date stored as literal
# text (e.g. "12/14/2025" in row 61), cols B/C/D are plain numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column A on this row to be treated as text so the date-like
# string isn't auto-converted into a date serial number, then drop the
# formatting Excel applies for that so the cell stays plain/unstyled.
$ws.Range("A62").NumberFormat = "@"
$ws.Range("A62").Value = "12/21/2025"
$ws.Range("A62").ClearFormats()

$ws.Range("B62").Value = 0.0005591900000000011
$ws.Range("C62").Value = 88520.89629642859
$ws.Range("D62").Value = 50
